# Update weekly triaged issues - shift months forward by one and update stats
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023 December", 55, 7),
    @("2024 January", 53, 7),
    @("2024 February", 54, 13),
    @("2024 March", 35, 7),
    @("2024 April", 36, 12),
    @("2024 May", 55, 20),
    @("2024 June", 41, 24),
    @("2024 July", 47, 18),
    @("2024 August", 35, 18),
    @("2024 September", 38, 19),
    @("2024 October", 29, 28),
    @("2024 November", 15, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    # Force text interpretation so "Month Year" strings are not
    # auto-converted into date serial numbers by Excel.
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
